$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted above the existing row 9, pushing
# the former rows 9-17 down to rows 10-18 (their data is unchanged).
$ws.Rows(9).Insert()

# Populate the newly inserted row 9 with the new record's data.
$ws.Range("A9").Value = 12
$ws.Range("B9").Value = "Mapocho Venta Directa de Santiago"
$ws.Range("C9").Value = "Metropolitana"
$ws.Range("D9").Value = 44581
$ws.Range("E9").Value = 13
$ws.Range("F9").Value = 100112021
$ws.Range("G9").Value = "Ají"
$ws.Range("H9").Value = "Americana (o)"
$ws.Range("I9").Value = "Segunda"
$ws.Range("J9").Value = 30
$ws.Range("K9").Value = 17000
$ws.Range("L9").Value = 17000
$ws.Range("M9").Value = 17000
$ws.Range("N9").Value = "`$/caja 25 kilos"
$ws.Range("O9").Value = "Provincia de Limarí"
$ws.Range("P9").Value = 680
$ws.Range("Q9").Value = 25
$ws.Range("R9").Value = "Hortaliza"
